$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '29.221.87'
$ws.Range('E2').Value = '  +0.62%  '

# Row 3
$ws.Range('D3').Value = '1.895.81'
$ws.Range('E3').Value = '  +0.37%  '

# Row 4
$ws.Range('D4').Value = '''1.001'
$ws.Range('E4').Value = '  +0.08%  '

# Row 5
$ws.Range('D5').Value = '''322.58'
$ws.Range('E5').Value = '  -2.45%  '

# Row 6
$ws.Range('E6').Value = '  +0.11%  '

# Row 7
$ws.Range('D7').Value = '''0.4699'

# Row 8
$ws.Range('D8').Value = '''0.4021'
$ws.Range('E8').Value = '  -1.94%  '

# Row 9
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').Value = '''0.07992'
$ws.Range('E9').Value = '  +0.56%  '

# Row 10
$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D10').Value = '''0.9941'
$ws.Range('E10').Value = '  -0.37%  '

# Row 11
$ws.Range('B11').Value = 'Solana'
$ws.Range('C11').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D11').Value = '''22.40'
$ws.Range('E11').Value = '  +2.53%  '

# Row 12
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.905.03'
$ws.Range('E12').Value = '  -0.71%  '

# Row 13
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '''5.853'
$ws.Range('E13').Value = '  -1.09%  '

# Row 14
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D14').Value = '''7.040'
$ws.Range('E14').Value = '  -0.54%  '

# Row 15
$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D15').Value = '''88.85'
$ws.Range('E15').Value = '  +0.18%  '

# Row 16
$ws.Range('D16').Value = '''1.001'
$ws.Range('E16').Value = '  +0.12%  '

# Row 17
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').Value = '''0.06608'
$ws.Range('E17').Value = '  +0.84%  '

# Row 18
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').Value = '''0.00001025'
$ws.Range('E18').Value = '  +0.38%  '

# Row 19
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D19').Value = '''17.43'
$ws.Range('E19').Value = '  -0.07%  '

# Row 20
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = '''1.002'
$ws.Range('E20').Value = '  +0.14%  '

# Row 21
$ws.Range('B21').Value = 'WrappedBTC'
$ws.Range('C21').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D21').Value = '29.243.84'
$ws.Range('E21').Value = '  +0.60%  '

# Row 22
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').Value = '''5.498'
$ws.Range('E22').Value = '  +1.08%  '

# Row 23
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D23').Value = '''11.59'
$ws.Range('E23').Value = '  +1.55%  '

# Row 24
$ws.Range('B24').Value = 'Toncoin'
$ws.Range('C24').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D24').Value = '''2.195'
$ws.Range('E24').Value = '  -0.19%  '

# Row 25
$ws.Range('B25').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C25').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D25').Value = '2.138.38'
$ws.Range('E25').Value = '  +0.07%  '

# Row 26
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '''154.89'
$ws.Range('E26').Value = '  -0.94%  '

# Row 27
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D27').Value = '''19.60'
$ws.Range('E27').Value = '  +0.23%  '

# Row 28
$ws.Range('B28').Value = 'InternetComputer(DFINITY)'
$ws.Range('C28').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D28').Value = '''6.057'
$ws.Range('E28').Value = '  +9.97%  '

# Row 29
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').Value = '''2.074'
$ws.Range('E29').Value = '  -0.74%  '

# Row 30
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = '''117.25'
$ws.Range('E30').Value = '  +0.04%  '

# Row 31
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').Value = '''1.046'
$ws.Range('E31').Value = '  +1.10%  '

# Row 32
$ws.Range('B32').Value = 'Stellar'
$ws.Range('C32').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D32').Value = '''0.09431'
$ws.Range('E32').Value = '  +1.13%  '

# Row 33
$ws.Range('B33').Value = 'ARBITRUM'
$ws.Range('C33').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D33').Value = '''1.392'
$ws.Range('E33').Value = '  -1.13%  '

# Row 34
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D34').Value = '''3.542'
$ws.Range('E34').Value = '  +0.66%  '

# Row 35
$ws.Range('B35').Value = 'Filecoin'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D35').Value = '''5.336'
$ws.Range('E35').Value = '  +0.88%  '

# Row 36
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').Value = '''0.06067'
$ws.Range('E36').Value = '  +0.25%  '

# Row 37
$ws.Range('B37').Value = 'VeChain'
$ws.Range('C37').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D37').Value = '''0.02232'
$ws.Range('E37').Value = '  +0.30%  '

# Row 38
$ws.Range('B38').Value = 'TrustWalletToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D38').Value = '''1.172'
$ws.Range('E38').Value = '  +0.02%  '

# Row 39
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').Value = '''8.057'
$ws.Range('E39').Value = '  -3.45%  '

# Row 40
$ws.Range('B40').Value = 'TheSandbox'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D40').Value = '''0.5806'
$ws.Range('E40').Value = '  +0.41%  '

# Row 41
$ws.Range('B41').Value = 'RenderToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D41').Value = '''2.479'
$ws.Range('E41').Value = '  +7.87%  '

# Row 42
$ws.Range('D42').Value = '''0.1822'
$ws.Range('E42').Value = '  +0.08%  '

# Row 43
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').Value = '''10.03'
$ws.Range('E43').Value = '  -0.60%  '

# Row 44
$ws.Range('B44').Value = 'WEMIXToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D44').Value = '''1.273'
$ws.Range('E44').Value = '  +1.16%  '

# Row 45
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').Value = '''0.07677'
$ws.Range('E45').Value = '  +2.45%  '

# Row 46
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D46').Value = '''12.09'
$ws.Range('E46').Value = '  +0.42%  '

# Row 47
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').Value = '''0.5460'
$ws.Range('E47').Value = '  +0.10%  '

# Row 48
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').Value = '''1.898'
$ws.Range('E48').Value = '  -0.32%  '

# Row 49
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').Value = '''113.53'
$ws.Range('E49').Value = '  +1.56%  '

# Row 50
$ws.Range('B50').Value = 'Elrond'
$ws.Range('C50').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D50').Value = '''44.06'
$ws.Range('E50').Value = '  +0.08%  '

# Row 51
$ws.Range('B51').Value = 'WOONetwork'
$ws.Range('C51').Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range('D51').Value = '''0.2931'
$ws.Range('E51').Value = '  +4.88%  '
